$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add thin bottom border under row 4 (new style xf 6 and 7) ---
$ws.Cells.Item(4,2).Borders.Item(9).LineStyle = 1
$ws.Cells.Item(4,1).Borders.Item(9).LineStyle = 1
$ws.Cells.Item(4,3).Borders.Item(9).LineStyle = 1
$ws.Cells.Item(4,4).Borders.Item(9).LineStyle = 1
$ws.Cells.Item(4,5).Borders.Item(9).LineStyle = 1

# --- Column C (English) rows 5-15 ---
$ws.Cells.Item(5,3).Value = ' Hey, you two.[K] Do you have any\nnews on [CS:P]Zero Isle[CR]?'
$ws.Cells.Item(6,3).Value = ' [CS:P]Zero Isle[CR]?'
$ws.Cells.Item(6,2).Value = '100, 104, 108'
$ws.Cells.Item(7,3).Value = ' Yeah. Rumor has it that the\ngreatest treasures are there…'
$ws.Cells.Item(8,3).Value = ' But the treasures aren\''t the\nmost fascinating thing about that place.'
$ws.Cells.Item(9,3).Value = ' When you go in the dungeon...[K]\neveryone is temporarily dropped to Level 1 and\ntheir learned abilities are zeroed out.'
$ws.Cells.Item(10,3).Value = ' From what we\''ve heard, you\ncan\''t even take in items.'
$ws.Cells.Item(11,3).Value = ' Basically, you walk in alone as if\nyou\''re a rank amateur again.'
$ws.Cells.Item(12,3).Value = ' No exploration team\''s made it to\nthe end yet.'
$ws.Cells.Item(13,3).Value = ' They\''re calling it the ultimate\ndungeon for good reason.'
$ws.Cells.Item(14,3).Value = ' [CS:N]Drapion[CR]![K] You talk too much!'
$ws.Cells.Item(15,3).Value = ' Urk![K] I\''m sorry, my lady [CS:N]Weavile[CR].'

# --- Column D (Russian) rows 5-15 ---
$ws.Cells.Item(5,4).Value = ' Эй, вы двое.[K] Вы знаете\nчто-нибудь о [CS:P]Нуль-Острове[CR]?'
$ws.Cells.Item(6,4).Value = ' [CS:P]Нуль-Острове[CR]?'
$ws.Cells.Item(7,4).Value = ' Да. Ходят слухи, что там есть\nсамые ценные сокровища...'
$ws.Cells.Item(8,4).Value = ' Но сокровища это не самое\nинтересное, что там есть.'
$ws.Cells.Item(9,4).Value = ' Если там пойти в подземелье...[K]\nУровень всей группы временно станет 1 и\nизученные способности будут обнулены.'
$ws.Cells.Item(10,4).Value = ' Из того, что известно, туда\nдаже нельзя взять предметы.'
$ws.Cells.Item(11,4).Value = ' По сути, вы снова становитесь\nзелёными новичками.'
$ws.Cells.Item(12,4).Value = ' Ещё ни одной команде\nисследователей не удалось покорить остров.'
$ws.Cells.Item(13,4).Value = ' Не зря эти подземелья называют\nтруднейшими.'
$ws.Cells.Item(14,4).Value = ' [CS:N]Драпион[CR]![K] Ты слишком много\nболтаешь!'
$ws.Cells.Item(15,4).Value = ' Ух![K] Прости, моя леди [CS:N]Вивайл[CR].'

# --- Column E (cipher) rows 5-15 ---
$ws.Cells.Item(5,5).Value = ' Üê, âú äâïå.[K] Âú èîàåóå\nœóï-îéáôäû ï [CS:P]Îôìû-Ïòóñïâå[CR]?'
$ws.Cells.Item(6,5).Value = ' [CS:P]Îôìû-Ïòóñïâå[CR]?'
$ws.Cells.Item(7,5).Value = ' Äà. Öïäÿó òìôöé, œóï óàí åòóû\nòàíúå øåîîúå òïëñïâéþà...'
$ws.Cells.Item(8,5).Value = ' Îï òïëñïâéþà üóï îå òàíïå\néîóåñåòîïå, œóï óàí åòóû.'
$ws.Cells.Item(9,5).Value = ' Åòìé óàí ðïêóé â ðïäèåíåìûå...[K]\nÔñïâåîû âòåê ãñôððú âñåíåîîï òóàîåó 1 é\néèôœåîîúå òðïòïáîïòóé áôäôó ïáîôìåîú.'
$ws.Cells.Item(10,5).Value = ' Éè óïãï, œóï éèâåòóîï, óôäà\näàçå îåìûèÿ âèÿóû ðñåäíåóú.'
$ws.Cells.Item(11,5).Value = ' Ðï òôóé, âú òîïâà òóàîïâéóåòû\nèåìæîúíé îïâéœëàíé.'
$ws.Cells.Item(12,5).Value = ' Åþæ îé ïäîïê ëïíàîäå\néòòìåäïâàóåìåê îå ôäàìïòû ðïëïñéóû ïòóñïâ.'
$ws.Cells.Item(13,5).Value = ' Îå èñÿ üóé ðïäèåíåìûÿ îàèúâàýó\nóñôäîåêšéíé.'
$ws.Cells.Item(14,5).Value = ' [CS:N]Äñàðéïî[CR]![K] Óú òìéšëïí íîïãï\náïìóàåšû!'
$ws.Cells.Item(15,5).Value = ' Ôö![K] Ðñïòóé, íïÿ ìåäé [CS:N]Âéâàêì[CR].'

# --- Column B numeric values ---
$ws.Cells.Item(5,2).Value = 85
$ws.Cells.Item(7,2).Value = 117
$ws.Cells.Item(8,2).Value = 120
$ws.Cells.Item(9,2).Value = 123
$ws.Cells.Item(10,2).Value = 126
$ws.Cells.Item(11,2).Value = 129
$ws.Cells.Item(12,2).Value = 132
$ws.Cells.Item(13,2).Value = 135
$ws.Cells.Item(14,2).Value = 144
$ws.Cells.Item(15,2).Value = 167

# --- Apply styles (wrap text, no border) to rows 5-15 like row 2/3 ---
$ws.Cells.Item(2,2).Copy()
$ws.Cells.Item(5,2).PasteSpecial(-4122)
$ws.Cells.Item(6,2).PasteSpecial(-4122)
$ws.Cells.Item(7,2).PasteSpecial(-4122)
$ws.Cells.Item(8,2).PasteSpecial(-4122)
$ws.Cells.Item(9,2).PasteSpecial(-4122)
$ws.Cells.Item(10,2).PasteSpecial(-4122)
$ws.Cells.Item(11,2).PasteSpecial(-4122)
$ws.Cells.Item(12,2).PasteSpecial(-4122)
$ws.Cells.Item(13,2).PasteSpecial(-4122)
$ws.Cells.Item(14,2).PasteSpecial(-4122)
$ws.Cells.Item(15,2).PasteSpecial(-4122)
$ws.Cells.Item(2,3).Copy()
$ws.Cells.Item(5,3).PasteSpecial(-4122)
$ws.Cells.Item(5,4).PasteSpecial(-4122)
$ws.Cells.Item(5,5).PasteSpecial(-4122)
$ws.Cells.Item(6,3).PasteSpecial(-4122)
$ws.Cells.Item(6,4).PasteSpecial(-4122)
$ws.Cells.Item(6,5).PasteSpecial(-4122)
$ws.Cells.Item(7,3).PasteSpecial(-4122)
$ws.Cells.Item(7,4).PasteSpecial(-4122)
$ws.Cells.Item(7,5).PasteSpecial(-4122)
$ws.Cells.Item(8,3).PasteSpecial(-4122)
$ws.Cells.Item(8,4).PasteSpecial(-4122)
$ws.Cells.Item(8,5).PasteSpecial(-4122)
$ws.Cells.Item(9,3).PasteSpecial(-4122)
$ws.Cells.Item(9,4).PasteSpecial(-4122)
$ws.Cells.Item(9,5).PasteSpecial(-4122)
$ws.Cells.Item(10,3).PasteSpecial(-4122)
$ws.Cells.Item(10,4).PasteSpecial(-4122)
$ws.Cells.Item(10,5).PasteSpecial(-4122)
$ws.Cells.Item(11,3).PasteSpecial(-4122)
$ws.Cells.Item(11,4).PasteSpecial(-4122)
$ws.Cells.Item(11,5).PasteSpecial(-4122)
$ws.Cells.Item(12,3).PasteSpecial(-4122)
$ws.Cells.Item(12,4).PasteSpecial(-4122)
$ws.Cells.Item(12,5).PasteSpecial(-4122)
$ws.Cells.Item(13,3).PasteSpecial(-4122)
$ws.Cells.Item(13,4).PasteSpecial(-4122)
$ws.Cells.Item(13,5).PasteSpecial(-4122)
$ws.Cells.Item(14,3).PasteSpecial(-4122)
$ws.Cells.Item(14,4).PasteSpecial(-4122)
$ws.Cells.Item(14,5).PasteSpecial(-4122)
$ws.Cells.Item(15,3).PasteSpecial(-4122)
$ws.Cells.Item(15,4).PasteSpecial(-4122)
$ws.Cells.Item(15,5).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row heights ---
$ws.Rows.Item(5).RowHeight = 21.6
$ws.Rows.Item(7).RowHeight = 21.6
$ws.Rows.Item(8).RowHeight = 21.6
$ws.Rows.Item(9).RowHeight = 52.2
$ws.Rows.Item(10).RowHeight = 21.6
$ws.Rows.Item(11).RowHeight = 21.6
$ws.Rows.Item(12).RowHeight = 21.6
$ws.Rows.Item(13).RowHeight = 21.6
$ws.Rows.Item(14).RowHeight = 21.6
$ws.Rows.Item(15).RowHeight = 21.6

# --- Sheet view adjustments ---
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("E15").Select()

